$d = $word.ActiveDocument

# 1. Update the iteration start date from 09 to 08.
$d.Content.Find.Execute("iteração durou entre as datas 09", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "iteração durou entre as datas 08", 2)

# 2. Move the "_GoBack" bookmark: remove it from its old location (end of the
#    document, after "...pesquisadas.") and re-create it right after the text
#    that was just edited ("...datas 08"), mirroring Word's behaviour of
#    tracking the most recent edit position with the hidden _GoBack bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute("iteração durou entre as datas 08", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$newBookmarkRange = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
